# Add team record (Wins / Losses / Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: copy formatting from the last existing header cell (AC1)
# so the new headers share the same bold/border/centered style (s="1"),
# then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows: every team record in this file is 65-97-0.
for ($r = 2; $r -le 51; $r++) {
    $ws.Range("AD$r").Value = 65
    $ws.Range("AE$r").Value = 97
    $ws.Range("AF$r").Value = 0
}
